$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at the top; this pushes all existing rows (and their
# styles/values) down by 3, matching the diff's row-shift (old row 1 -> 4,
# old row 2 -> 5, ..., old row 22 -> 25).
$ws.Rows("1:3").Insert()

# Populate the three new description rows in column A (unstyled / default style).
$ws.Range("A1").Value = "This is the original case study area MAR."
$ws.Range("A2").Value = "It includes quantitative and qualitative groundwater data for 10 monitoring sites."
$ws.Range("A3").Value = "There are 20 different parameters measured. In total, there are 830 groundwater samples with over 1,400 individual datapoints."
